$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1514.5834
$ws.Range("I2").Value = 171
$ws.Range("K2").Value = 171
$ws.Range("M2").Value = -58

$ws.Range("H17").Value = 1486.5
$ws.Range("J17").Value = 1486.5
$ws.Range("L17").Value = 4459.5
$ws.Range("N17").Value = -4795.5

$ws.Range("H38").Value = 2696.5881
$ws.Range("I38").Value = 209.625
$ws.Range("K38").Value = 628.875
$ws.Range("M38").Value = -256.875

$ws.Range("H42").Value = 1028.625
$ws.Range("I42").Value = 148.16667
$ws.Range("J42").Value = 1322.1111
$ws.Range("K42").Value = 444.50001
$ws.Range("L42").Value = 3966.3333
$ws.Range("M42").Value = -214.50001
$ws.Range("N42").Value = -4426.3333

$ws.Range("H43").Value = 6534.75
$ws.Range("I43").Value = 6741.4
$ws.Range("J43").Value = 6190.3335
$ws.Range("K43").Value = 6741.4
$ws.Range("L43").Value = 6190.3335
$ws.Range("M43").Value = -6672.4
$ws.Range("N43").Value = -6328.3335

$ws.Range("H46").Value = 3647.25
$ws.Range("I46").Value = 3295
$ws.Range("J46").Value = 3999.5
$ws.Range("K46").Value = 9885
$ws.Range("L46").Value = 11998.5
$ws.Range("M46").Value = -9766
$ws.Range("N46").Value = -12236.5

$ws.Range("H60").Value = 3647.25
$ws.Range("I60").Value = 3295
$ws.Range("J60").Value = 3999.5
$ws.Range("K60").Value = 9885
$ws.Range("L60").Value = 11998.5
$ws.Range("M60").Value = -9401
$ws.Range("N60").Value = -12966.5

$ws.Range("H74").Value = 6549.4287
$ws.Range("I74").Value = 7041
$ws.Range("K74").Value = 7041
$ws.Range("M74").Value = -6105

$ws.Range("H76").Value = 5086.923
$ws.Range("I76").Value = 3209
$ws.Range("J76").Value = 6260.625
$ws.Range("K76").Value = 3209
$ws.Range("L76").Value = 6260.625
$ws.Range("M76").Value = -2894
$ws.Range("N76").Value = -6890.625

$ws.Range("H77").Value = 6549.4287
$ws.Range("I77").Value = 7041
$ws.Range("K77").Value = 35205
$ws.Range("M77").Value = -30525

$ws.Range("H79").Value = 5086.923
$ws.Range("I79").Value = 3209
$ws.Range("J79").Value = 6260.625
$ws.Range("K79").Value = 3209
$ws.Range("L79").Value = 6260.625
$ws.Range("M79").Value = -2117
$ws.Range("N79").Value = -8444.625

$ws.Range("H128").Value = 90000
$ws.Range("J128").Value = 90000
$ws.Range("L128").Value = 90000
$ws.Range("N128").Value = -99960

$ws.Range("H137").Value = 2593.2334
$ws.Range("I137").Value = 1418.125
$ws.Range("K137").Value = 4254.375
$ws.Range("M137").Value = -1704.375

$ws.Range("H138").Value = 3280.9062
$ws.Range("I138").Value = 2400.3462
$ws.Range("K138").Value = 7201.0386
$ws.Range("M138").Value = -2061.0386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 136.42857
$ws.Range("I5").Value = 153.83333
$ws.Range("K5").Value = 153.83333
$ws.Range("M5").Value = -41.83332999999999

$ws.Range("H37").Value = 24499.5
$ws.Range("I37").Value = 24499.5
$ws.Range("K37").Value = 24499.5
$ws.Range("M37").Value = -24226.5

$ws.Range("H61").Value = 10493.392
$ws.Range("I61").Value = 8245.875
$ws.Range("K61").Value = 8245.875
$ws.Range("M61").Value = -8033.875

$ws.Range("H97").Value = 1630.7778
$ws.Range("I97").Value = 1702.875
$ws.Range("K97").Value = 1702.875
$ws.Range("M97").Value = -1206.875

$ws.Range("H132").Value = 4347.115
$ws.Range("I132").Value = 1967.4
$ws.Range("K132").Value = 5902.200000000001
$ws.Range("M132").Value = -3372.200000000001

$ws.Range("H136").Value = 10493.392
$ws.Range("I136").Value = 8245.875
$ws.Range("K136").Value = 24737.625
$ws.Range("M136").Value = -22187.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 136.42857
$ws.Range("I4").Value = 153.83333
$ws.Range("K4").Value = 153.83333
$ws.Range("M4").Value = -38.83332999999999

$ws.Range("H20").Value = 1324.3846
$ws.Range("I20").Value = 1466.5238
$ws.Range("J20").Value = 727.4
$ws.Range("K20").Value = 1466.5238
$ws.Range("L20").Value = 727.4
$ws.Range("M20").Value = -1219.5238
$ws.Range("N20").Value = -1221.4

$ws.Range("H22").Value = 510.13635
$ws.Range("I22").Value = 479.33334
$ws.Range("K22").Value = 479.33334
$ws.Range("M22").Value = -306.33334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6516
$ws.Range("I31").Value = 2745.4707
$ws.Range("K31").Value = 2745.4707
$ws.Range("M31").Value = -2450.4707

$ws.Range("H34").Value = 6516
$ws.Range("I34").Value = 2745.4707
$ws.Range("K34").Value = 2745.4707
$ws.Range("M34").Value = -2543.4707

$ws.Range("H134").Value = 6853.5454
$ws.Range("I134").Value = 3199.6
$ws.Range("K134").Value = 9598.799999999999
$ws.Range("M134").Value = -7063.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2113.6667
$ws.Range("I5").Value = 2436.4
$ws.Range("K5").Value = 7309.200000000001
$ws.Range("M5").Value = -7197.200000000001

$ws.Range("H122").Value = 870.6
$ws.Range("J122").Value = 884.5
$ws.Range("L122").Value = 7960.5
$ws.Range("N122").Value = -12860.5

$ws.Range("H135").Value = 2113.6667
$ws.Range("I135").Value = 2436.4
$ws.Range("K135").Value = 21927.6
$ws.Range("M135").Value = -19392.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 41599.4
$ws.Range("J43").Value = 49333
$ws.Range("L43").Value = 49333
$ws.Range("N43").Value = -49635

$ws.Range("H46").Value = 29343.666
$ws.Range("I46").Value = 9020.5
$ws.Range("K46").Value = 9020.5
$ws.Range("M46").Value = -8864.5

$ws.Range("H57").Value = 39750
$ws.Range("I57").Value = 35000
$ws.Range("J57").Value = 54000
$ws.Range("K57").Value = 35000
$ws.Range("L57").Value = 54000
$ws.Range("M57").Value = -34180
$ws.Range("N57").Value = -55640

$ws.Range("H70").Value = 3899.7612
$ws.Range("I70").Value = 1997.4546
$ws.Range("J70").Value = 4273.4287
$ws.Range("K70").Value = 1997.4546
$ws.Range("L70").Value = 4273.4287
$ws.Range("M70").Value = -1727.4546
$ws.Range("N70").Value = -4813.4287

$ws.Range("H73").Value = 3899.7612
$ws.Range("I73").Value = 1997.4546
$ws.Range("J73").Value = 4273.4287
$ws.Range("K73").Value = 1997.4546
$ws.Range("L73").Value = 4273.4287
$ws.Range("M73").Value = -1061.4546
$ws.Range("N73").Value = -6145.4287

$ws.Range("H97").Value = 566.9
$ws.Range("I97").Value = 628.1667
$ws.Range("J97").Value = 475
$ws.Range("K97").Value = 628.1667
$ws.Range("L97").Value = 475
$ws.Range("M97").Value = -132.1667
$ws.Range("N97").Value = -1467

$ws.Range("H126").Value = 3600.8
$ws.Range("I126").Value = 2667.6667
$ws.Range("K126").Value = 8003.000100000001
$ws.Range("M126").Value = -5533.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1664.4667
$ws.Range("I16").Value = 1865.8334
$ws.Range("K16").Value = 1865.8334
$ws.Range("M16").Value = -1695.8334

$ws.Range("H55").Value = 1280.5714
$ws.Range("I55").Value = 629.625
$ws.Range("K55").Value = 629.625
$ws.Range("M55").Value = -456.625

$ws.Range("H82").Value = 1001.5238
$ws.Range("I82").Value = 652.2222
$ws.Range("K82").Value = 652.2222
$ws.Range("M82").Value = -291.2222

$ws.Range("H85").Value = 1001.5238
$ws.Range("I85").Value = 652.2222
$ws.Range("K85").Value = 652.2222
$ws.Range("M85").Value = 595.7778

$ws.Range("H87").Value = 10000
$ws.Range("I87").Value = 10000
$ws.Range("K87").Value = 10000
$ws.Range("M87").Value = -8877

$ws.Range("H90").Value = 10000
$ws.Range("I90").Value = 10000
$ws.Range("K90").Value = 30000
$ws.Range("M90").Value = -24384

$ws.Range("H100").Value = 5011
$ws.Range("I100").Value = 4579.1113
$ws.Range("J100").Value = 5982.75
$ws.Range("K100").Value = 4579.1113
$ws.Range("L100").Value = 5982.75
$ws.Range("M100").Value = -4038.1113
$ws.Range("N100").Value = -7064.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2454.6875
$ws.Range("J81").Value = 2663
$ws.Range("L81").Value = 5326
$ws.Range("N81").Value = -7448

$ws.Range("H84").Value = 2454.6875
$ws.Range("J84").Value = 2663
$ws.Range("L84").Value = 26630
$ws.Range("N84").Value = -37238

$ws.Range("H141").Value = 74317
$ws.Range("J141").Value = 74317
$ws.Range("L141").Value = 74317
$ws.Range("N141").Value = -84677
